# Insert two new data rows (23 and 24) into the "Membrillo" sheet, pushing the
# existing rows 23..119 down to 25..121, then populate the two new rows with
# fresh observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 23:24 downward (and everything below them) by inserting two
# blank rows at position 23.
$ws.Range("A23:T24").EntireRow.Insert()

# --- New row 23 -----------------------------------------------------------
$ws.Range("A23").Value = 6
$ws.Range("B23").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C23").Value = "Metropolitana"
$ws.Range("D23").Value = 44715
$ws.Range("E23").Value = 13
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100104
$ws.Range("H23").Value = "Frutos de pepita"
$ws.Range("I23").Value = 100104003
$ws.Range("J23").Value = "Membrillo"
$ws.Range("K23").Value = "Champion"
$ws.Range("L23").Value = "Especial"
$ws.Range("M23").Value = 5
$ws.Range("N23").Value = 270000
$ws.Range("O23").Value = 270000
$ws.Range("P23").Value = 270000
$ws.Range("Q23").Value = "$/bins (450 kilos)"
$ws.Range("R23").Value = "Región de O'Higgins"
$ws.Range("S23").Value = 600
$ws.Range("T23").Value = 450

# --- New row 24 -----------------------------------------------------------
$ws.Range("A24").Value = 6
$ws.Range("B24").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C24").Value = "Metropolitana"
$ws.Range("D24").Value = 44715
$ws.Range("E24").Value = 13
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100104
$ws.Range("H24").Value = "Frutos de pepita"
$ws.Range("I24").Value = 100104003
$ws.Range("J24").Value = "Membrillo"
$ws.Range("K24").Value = "Champion"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 8
$ws.Range("N24").Value = 200000
$ws.Range("O24").Value = 200000
$ws.Range("P24").Value = 200000
$ws.Range("Q24").Value = "$/bins (450 kilos)"
$ws.Range("R24").Value = "Región de O'Higgins"
$ws.Range("S24").Value = 444
$ws.Range("T24").Value = 450
